# Append the new resale-number row (row 60) reported 2024-01-15 18:32:10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 60

# Columns A, B, D hold text that looks like a date/time/number to Excel's
# auto-conversion ("2024-01-15", "18:32:10", "02"). Force them to stay text
# (matching the existing rows, which are plain inline strings with no
# special formatting) by temporarily marking the cell as Text before the
# write, then re-syncing the cell style back to an unformatted neighbour so
# no stray number-format style is left behind on the new cells.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-15"
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item(1, 3).Style

$ws.Cells.Item($row, 2).Value = "18:32:10"

$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item(1, 3).Style

$ws.Cells.Item($row, 5).Value = 138902
$ws.Cells.Item($row, 6).Value = 139087
$ws.Cells.Item($row, 7).Value = 171397
$ws.Cells.Item($row, 8).Value = 148148
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119390
$ws.Cells.Item($row, 11).Value = 220952
$ws.Cells.Item($row, 12).Value = 254084
$ws.Cells.Item($row, 13).Value = 184962
$ws.Cells.Item($row, 14).Value = 110382
$ws.Cells.Item($row, 15).Value = 41130
$ws.Cells.Item($row, 16).Value = 30885
$ws.Cells.Item($row, 17).Value = 73186
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42560
$ws.Cells.Item($row, 20).Value = -1
